$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet data (rows 2-9) is being fully replaced by a reordered/updated
# table spanning rows 2-12 (3 new observation rows added). Clear the old
# data first, then rewrite every row from scratch.
$ws.Range("A2:AY9").ClearContents()

# Startdatum/Slutdatum (Y, AA) hold literal date-like text (e.g. "2023-09-13"),
# not real Excel dates. Force Text number format first so the assignment below
# is not reinterpreted as a date serial.
$ws.Range("Y2:Y12").NumberFormat = "@"
$ws.Range("AA2:AA12").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 112067161
$ws.Range("B2").Value = 88167
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 1599
$ws.Range("F2").Value = "Fjällfotad musseron"
$ws.Range("G2").Value = "Tricholoma olivaceotinctum"
$ws.Range("H2").Value = "Mort.Chr. & Heilm.-Claus."
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("N2").Value = ""
$ws.Range("P2").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q2").Value = 469256
$ws.Range("R2").Value = 7039724
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Jämtland"
$ws.Range("U2").Value = "Krokom"
$ws.Range("V2").Value = "Jämtland"
$ws.Range("W2").Value = "Offerdal"
$ws.Range("Y2").Value = "2023-09-13"
$ws.Range("AA2").Value = "2023-09-13"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AF2").Value = ""
$ws.Range("AG2").Value = $false
$ws.Range("AT2").Value = ""
$ws.Range("AW2").Value = "Rashid Kadhim"
$ws.Range("AX2").Value = "Rashid Kadhim"
$ws.Range("AY2").Value = ""

# Row 3
$ws.Range("A3").Value = 112074141
$ws.Range("B3").Value = 90434
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 4745
$ws.Range("F3").Value = "Tallriska"
$ws.Range("G3").Value = "Lactarius musteus"
$ws.Range("H3").Value = "Fr."
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("P3").Value = "Tjärnmyren, Nybodarna, Offerdal, Jmt"
$ws.Range("Q3").Value = 469230
$ws.Range("R3").Value = 7039721
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Jämtland"
$ws.Range("U3").Value = "Krokom"
$ws.Range("V3").Value = "Jämtland"
$ws.Range("W3").Value = "Offerdal"
$ws.Range("Y3").Value = "2023-09-13"
$ws.Range("AA3").Value = "2023-09-13"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AF3").Value = ""
$ws.Range("AG3").Value = $false
$ws.Range("AT3").Value = ""
$ws.Range("AW3").Value = "Rashid Kadhim"
$ws.Range("AX3").Value = "Rashid Kadhim"
$ws.Range("AY3").Value = ""

# Row 4
$ws.Range("A4").Value = 112370021
$ws.Range("B4").Value = 56430
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("I4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = ""
$ws.Range("P4").Value = "nybodarna Österulvsås, Jmt"
$ws.Range("Q4").Value = 469287
$ws.Range("R4").Value = 7039645
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Jämtland"
$ws.Range("U4").Value = "Krokom"
$ws.Range("V4").Value = "Jämtland"
$ws.Range("W4").Value = "Offerdal"
$ws.Range("Y4").Value = "2023-09-26"
$ws.Range("AA4").Value = "2023-09-26"
$ws.Range("AC4").Value = "ringhack"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Benny Öwre"
$ws.Range("AX4").Value = "Benny Öwre"
$ws.Range("AY4").Value = ""

# Row 5
$ws.Range("A5").Value = 112370020
$ws.Range("B5").Value = 56430
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("I5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "nybodarna Österulvsås, Jmt"
$ws.Range("Q5").Value = 469262
$ws.Range("R5").Value = 7039652
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Krokom"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Offerdal"
$ws.Range("Y5").Value = "2023-09-26"
$ws.Range("AA5").Value = "2023-09-26"
$ws.Range("AC5").Value = "ringhack äldre"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Value = ""
$ws.Range("AW5").Value = "Benny Öwre"
$ws.Range("AX5").Value = "Benny Öwre"
$ws.Range("AY5").Value = ""

# Row 6
$ws.Range("A6").Value = 112067953
$ws.Range("B6").Value = 88166
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "VU"
$ws.Range("E6").Value = 6276
$ws.Range("F6").Value = "Goliatmusseron"
$ws.Range("G6").Value = "Tricholoma matsutake"
$ws.Range("H6").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("P6").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q6").Value = 469442
$ws.Range("R6").Value = 7039562
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Krokom"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Offerdal"
$ws.Range("Y6").Value = "2023-09-13"
$ws.Range("AA6").Value = "2023-09-13"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Value = ""
$ws.Range("AW6").Value = "Rashid Kadhim"
$ws.Range("AX6").Value = "Rashid Kadhim"
$ws.Range("AY6").Value = ""

# Row 7
$ws.Range("A7").Value = 112068040
$ws.Range("B7").Value = 90816
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 2059
$ws.Range("F7").Value = "Skrovlig taggsvamp"
$ws.Range("G7").Value = "Hydnellum scabrosum"
$ws.Range("H7").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("P7").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q7").Value = 469465
$ws.Range("R7").Value = 7039571
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Jämtland"
$ws.Range("U7").Value = "Krokom"
$ws.Range("V7").Value = "Jämtland"
$ws.Range("W7").Value = "Offerdal"
$ws.Range("Y7").Value = "2023-09-13"
$ws.Range("AA7").Value = "2023-09-13"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = ""
$ws.Range("AW7").Value = "Rashid Kadhim"
$ws.Range("AX7").Value = "Rashid Kadhim"
$ws.Range("AY7").Value = ""

# Row 8
$ws.Range("A8").Value = 112068136
$ws.Range("B8").Value = 88166
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 6276
$ws.Range("F8").Value = "Goliatmusseron"
$ws.Range("G8").Value = "Tricholoma matsutake"
$ws.Range("H8").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("P8").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q8").Value = 469497
$ws.Range("R8").Value = 7039592
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Jämtland"
$ws.Range("U8").Value = "Krokom"
$ws.Range("V8").Value = "Jämtland"
$ws.Range("W8").Value = "Offerdal"
$ws.Range("Y8").Value = "2023-09-13"
$ws.Range("AA8").Value = "2023-09-13"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = "Rashid Kadhim"
$ws.Range("AX8").Value = "Rashid Kadhim"
$ws.Range("AY8").Value = ""

# Row 9
$ws.Range("A9").Value = 112067971
$ws.Range("B9").Value = 90785
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1968
$ws.Range("F9").Value = "Grantaggsvamp"
$ws.Range("G9").Value = "Bankera violascens"
$ws.Range("H9").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("I9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("P9").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q9").Value = 469442
$ws.Range("R9").Value = 7039562
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Jämtland"
$ws.Range("U9").Value = "Krokom"
$ws.Range("V9").Value = "Jämtland"
$ws.Range("W9").Value = "Offerdal"
$ws.Range("Y9").Value = "2023-09-13"
$ws.Range("AA9").Value = "2023-09-13"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AT9").Value = ""
$ws.Range("AW9").Value = "Rashid Kadhim"
$ws.Range("AX9").Value = "Rashid Kadhim"
$ws.Range("AY9").Value = ""

# Row 10
$ws.Range("A10").Value = 112068010
$ws.Range("B10").Value = 88166
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 6276
$ws.Range("F10").Value = "Goliatmusseron"
$ws.Range("G10").Value = "Tricholoma matsutake"
$ws.Range("H10").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("P10").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q10").Value = 469452
$ws.Range("R10").Value = 7039595
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = "Jämtland"
$ws.Range("U10").Value = "Krokom"
$ws.Range("V10").Value = "Jämtland"
$ws.Range("W10").Value = "Offerdal"
$ws.Range("Y10").Value = "2023-09-13"
$ws.Range("AA10").Value = "2023-09-13"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AT10").Value = ""
$ws.Range("AW10").Value = "Rashid Kadhim"
$ws.Range("AX10").Value = "Rashid Kadhim"
$ws.Range("AY10").Value = ""

# Row 11
$ws.Range("A11").Value = 112068431
$ws.Range("B11").Value = 90816
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 2059
$ws.Range("F11").Value = "Skrovlig taggsvamp"
$ws.Range("G11").Value = "Hydnellum scabrosum"
$ws.Range("H11").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("N11").Value = ""
$ws.Range("P11").Value = "Tjärnmyren (Tjärnmyren), Jmt"
$ws.Range("Q11").Value = 469560
$ws.Range("R11").Value = 7039585
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = "Jämtland"
$ws.Range("U11").Value = "Krokom"
$ws.Range("V11").Value = "Jämtland"
$ws.Range("W11").Value = "Offerdal"
$ws.Range("Y11").Value = "2023-09-13"
$ws.Range("AA11").Value = "2023-09-13"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AF11").Value = ""
$ws.Range("AG11").Value = $false
$ws.Range("AT11").Value = ""
$ws.Range("AW11").Value = "Rashid Kadhim"
$ws.Range("AX11").Value = "Rashid Kadhim"
$ws.Range("AY11").Value = ""

# Row 12
$ws.Range("A12").Value = 112370012
$ws.Range("B12").Value = 56430
$ws.Range("C12").Value = "Ovaliderad"
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("I12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("L12").Value = ""
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("P12").Value = "nybodarna Österulvsås, Jmt"
$ws.Range("Q12").Value = 469631
$ws.Range("R12").Value = 7039391
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = "Jämtland"
$ws.Range("U12").Value = "Krokom"
$ws.Range("V12").Value = "Jämtland"
$ws.Range("W12").Value = "Offerdal"
$ws.Range("Y12").Value = "2023-09-28"
$ws.Range("AA12").Value = "2023-09-28"
$ws.Range("AC12").Value = "ringhack"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AT12").Value = ""
$ws.Range("AW12").Value = "Benny Öwre"
$ws.Range("AX12").Value = "Benny Öwre"
$ws.Range("AY12").Value = ""

# Restore the default (General) style on the date-text columns now that the
# literal values are safely stored as text, so no stray number format lingers.
$ws.Range("Y2:Y12").Style = "Normal"
$ws.Range("AA2:AA12").Style = "Normal"
